# Update the cryptos price/volume table (columns D and E, rows 2-51) with
# freshly scraped values, as produced by the GitHub Actions refresh job.
#
# Column D ("Price") values are leading-apostrophe-prefixed so Excel keeps
# them as literal text (matching the sheet's existing inline-string cells)
# instead of auto-coercing them into numbers and dropping formatting such
# as trailing zeros (e.g. "8.80" -> 8.8) or thousand-dot grouping
# (e.g. "69.062.70").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.062.70"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "'2.499.72"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'570.33"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").Value = "'164.88"
$ws.Range("E6").Value = "  -2.82%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.511"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "'2.497.61"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").Value = "'0.159"
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("D11").Value = "'0.167"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "'0.354"
$ws.Range("E12").Value = "  +2.95%  "
$ws.Range("D13").Value = "'4.88"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").Value = "'2.963.78"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").Value = "'69.036.95"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "'0.0000174"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").Value = "'24.64"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").Value = "'2.503.38"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'11.22"
$ws.Range("E19").Value = "  -2.62%  "
$ws.Range("D20").Value = "'7.62"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "'346.31"
$ws.Range("E21").Value = "  -2.51%  "
$ws.Range("D22").Value = "'3.88"
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("D23").Value = "'1.98"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D25").Value = "'70.03"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").Value = "'3.89"
$ws.Range("E26").Value = "  -4.74%  "
$ws.Range("D27").Value = "'8.80"
$ws.Range("E27").Value = "  -4.39%  "
$ws.Range("D28").Value = "'2.649.15"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "'0.0₃0880"
$ws.Range("E30").Value = "  -3.32%  "
$ws.Range("D31").Value = "'7.74"
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("D32").Value = "'457.15"
$ws.Range("E32").Value = "  -5.53%  "
$ws.Range("D33").Value = "'1.22"
$ws.Range("E33").Value = "  -5.66%  "
$ws.Range("D34").Value = "'1.72"
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "'0.115"
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("D37").Value = "'155.57"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").Value = "'19.00"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").Value = "'18.37"
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'0.316"
$ws.Range("E41").Value = "  -1.58%  "
$ws.Range("D42").Value = "'4.65"
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("D43").Value = "'1.59"
$ws.Range("E43").Value = "  -3.60%  "
$ws.Range("D44").Value = "'38.06"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("E45").Value = "  -9.04%  "
$ws.Range("D46").Value = "'2.20"
$ws.Range("E46").Value = "  -7.95%  "
$ws.Range("D47").Value = "'141.23"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("D48").Value = "'0.518"
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("D49").Value = "'3.44"
$ws.Range("E49").Value = "  -2.28%  "
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").Value = "'0.574"
$ws.Range("E51").Value = "  -4.21%  "
